$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = 123
$ws.Range("B2").Value = 234
$ws.Range("C3").Value = 345
$ws.Range("D4").Value = 456
$ws.Range("E5").Value = 567

$ws.Range("E5").Select()
